$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price strings in column D stay as text (matching the
# original inline-string cell type) instead of being auto-coerced to numbers.
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '56.332.85'
$ws.Range('E2').Value = '  -3.16%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.966.68'
$ws.Range('E3').Value = '  -5.35%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '494.61'
$ws.Range('E5').Value = '  -5.63%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '134.30'
$ws.Range('E6').Value = '  +0.13%  '
$ws.Range('E7').Value = '  -0.04%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '2.969.03'
$ws.Range('E8').Value = '  -5.23%  '
$ws.Range('E9').Value = '  -4.15%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '7.16'
$ws.Range('E10').Value = '  -1.05%  '
$ws.Range('E11').Value = '  -3.76%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.349'
$ws.Range('E12').Value = '  -7.40%  '
$ws.Range('E13').Value = '  -0.78%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '3.478.61'
$ws.Range('E14').Value = '  -5.25%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '25.02'
$ws.Range('E15').Value = '  -1.57%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '56.328.34'
$ws.Range('E16').Value = '  -3.17%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.972.06'
$ws.Range('E17').Value = '  -5.23%  '
$ws.Range('E18').Value = '  -4.55%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '5.76'
$ws.Range('E19').Value = '  +0.51%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '12.32'
$ws.Range('E20').Value = '  -5.20%  '
$ws.Range('E21').Value = '  -2.00%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '324.86'
$ws.Range('E22').Value = '  -5.28%  '
$ws.Range('E23').Value = '  -0.11%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.466'
$ws.Range('E24').Value = '  -8.05%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '61.39'
$ws.Range('E25').Value = '  -9.41%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.00'
$ws.Range('E26').Value = '  +0.92%  '
$ws.Range('E27').Value = '  -5.96%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.0₃0893'
$ws.Range('E28').Value = '  -6.22%  '
$ws.Range('E29').Value = '  -0.06%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '6.48'
$ws.Range('E30').Value = '  -4.80%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '6.70'
$ws.Range('E31').Value = '  -2.83%  '
$ws.Range('E32').Value = '  -5.76%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.73'
$ws.Range('E33').Value = '  -6.73%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '20.14'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '152.75'
$ws.Range('E35').Value = '  -2.92%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '4.42'
$ws.Range('E36').Value = '  -8.52%  '
$ws.Range('E37').Value = '  -7.04%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '5.58'
$ws.Range('E38').Value = '  -10.29%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.0667'
$ws.Range('E39').Value = '  -2.59%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '23.09'
$ws.Range('E40').Value = '  -3.27%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '3.003.53'
$ws.Range('E41').Value = '  -5.11%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '36.57'
$ws.Range('E42').Value = '  -9.60%  '
$ws.Range('E43').Value = '  +0.02%  '
$ws.Range('E44').Value = '  -7.85%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.993'
$ws.Range('E45').Value = '  -9.06%  '
$ws.Range('B46').Value = 'Maker'
$ws.Range('C46').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.202.69'
$ws.Range('E46').Value = '  -3.55%  '
$ws.Range('B47').Value = 'Stacks'
$ws.Range('C47').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.39'
$ws.Range('E47').Value = '  -3.46%  '
$ws.Range('E48').Value = '  -9.12%  '
$ws.Range('E49').Value = '  +3.52%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0236'
$ws.Range('E50').Value = '  +0.88%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '5.72'
$ws.Range('E51').Value = '  -7.03%  '
